$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: 04-08-2021
# Write the date label as a text formula first, then paste its computed
# value back over itself as values-only, so Excel stores it as plain text
# (matching the existing "DD-MM-YYYY" labels in column A) instead of
# re-parsing the literal and auto-converting it to a date serial.
$ws.Range("A6").Formula = '="04-08-2021"'
$ws.Range("A6").Copy() | Out-Null
$ws.Range("A6").PasteSpecial(-4163) | Out-Null
$ws.Range("B6").Value = 10000
$ws.Range("D6").Value = 0

# Row 7: 05-08-2021
$ws.Range("A7").Formula = '="05-08-2021"'
$ws.Range("A7").Copy() | Out-Null
$ws.Range("A7").PasteSpecial(-4163) | Out-Null
$ws.Range("B7").Value = 10000
$ws.Range("C7").Value = 8000
$ws.Range("D7").Value = 8000
$ws.Range("E7").Value = 7000
$ws.Range("F7").Value = 1000
$ws.Range("G7").Value = 2.05

$excel.CutCopyMode = 0
